$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 47.86240033333333
$ws.Cells.Item(2, 8).Value = 143.587201
$ws.Cells.Item(2, 9).Value = 0.1228118231805696
$ws.Cells.Item(2, 10).Value = 0.1228118231805696
$ws.Cells.Item(2, 13).Value = 14.13421233333333
$ws.Cells.Item(2, 14).Value = 42.402637
$ws.Cells.Item(2, 15).Value = 0.06429893302526193
$ws.Cells.Item(2, 16).Value = 0.06429893302526192
$ws.Cells.Item(2, 17).Value = 676.4973290943374
$ws.Cells.Item(2, 18).Value = 6088.475961849036
$ws.Cells.Item(2, 19).Value = 0.007896669193397758
$ws.Cells.Item(2, 20).Value = 0.007896669193397755

$ws.Cells.Item(3, 7).Value = 47.86240033333333
$ws.Cells.Item(3, 8).Value = 143.587201
$ws.Cells.Item(3, 9).Value = 0.1228118231805696
$ws.Cells.Item(3, 10).Value = 0.1228118231805696
$ws.Cells.Item(3, 14).Value = 49.84314599999999
$ws.Cells.Item(3, 15).Value = 0.07558164617031606
$ws.Cells.Item(3, 16).Value = 0.07558164617031604
$ws.Cells.Item(3, 17).Value = 795.2042025749271
$ws.Cells.Item(3, 18).Value = 7156.837823174345
$ws.Cells.Item(3, 19).Value = 0.009282319765165232
$ws.Cells.Item(3, 20).Value = 0.00928231976516523

$ws.Cells.Item(4, 7).Value = 47.86240033333333
$ws.Cells.Item(4, 8).Value = 143.587201
$ws.Cells.Item(4, 9).Value = 0.1228118231805696
$ws.Cells.Item(4, 10).Value = 0.1228118231805696
$ws.Cells.Item(4, 13).Value = 87.857732
$ws.Cells.Item(4, 14).Value = 263.573196
$ws.Cells.Item(4, 15).Value = 0.3996797481453391
$ws.Cells.Item(4, 16).Value = 0.399679748145339
$ws.Cells.Item(4, 17).Value = 4205.08194136271
$ws.Cells.Item(4, 18).Value = 37845.73747226439
$ws.Cells.Item(4, 19).Value = 0.04908539855807999
$ws.Cells.Item(4, 20).Value = 0.04908539855807997

$ws.Cells.Item(5, 7).Value = 47.86240033333333
$ws.Cells.Item(5, 8).Value = 143.587201
$ws.Cells.Item(5, 9).Value = 0.1228118231805696
$ws.Cells.Item(5, 10).Value = 0.1228118231805696
$ws.Cells.Item(5, 13).Value = 5.425038666666667
$ws.Cells.Item(5, 14).Value = 16.275116
$ws.Cells.Item(5, 15).Value = 0.02467942250059516
$ws.Cells.Item(5, 16).Value = 0.02467942250059516
$ws.Cells.Item(5, 17).Value = 259.6553724878129
$ws.Cells.Item(5, 18).Value = 2336.898352390316
$ws.Cells.Item(5, 19).Value = 0.003030924872341665
$ws.Cells.Item(5, 20).Value = 0.003030924872341664

$ws.Cells.Item(6, 7).Value = 47.86240033333333
$ws.Cells.Item(6, 8).Value = 143.587201
$ws.Cells.Item(6, 9).Value = 0.1228118231805696
$ws.Cells.Item(6, 10).Value = 0.1228118231805696
$ws.Cells.Item(6, 13).Value = 95.78895966666666
$ws.Cells.Item(6, 14).Value = 287.366879
$ws.Cells.Item(6, 15).Value = 0.4357602501584878
$ws.Cells.Item(6, 16).Value = 0.4357602501584877
$ws.Cells.Item(6, 17).Value = 4584.689535079519
$ws.Cells.Item(6, 18).Value = 41262.20581571567
$ws.Cells.Item(6, 19).Value = 0.05351651079158499
$ws.Cells.Item(6, 20).Value = 0.05351651079158497

$ws.Cells.Item(7, 9).Value = 0.04786922362394307
$ws.Cells.Item(7, 10).Value = 0.04786922362394307
$ws.Cells.Item(7, 13).Value = 14.13421233333333
$ws.Cells.Item(7, 14).Value = 42.402637
$ws.Cells.Item(7, 15).Value = 0.06429893302526193
$ws.Cells.Item(7, 16).Value = 0.06429893302526192
$ws.Cells.Item(7, 17).Value = 263.6830973497062
$ws.Cells.Item(7, 18).Value = 2373.147876147356
$ws.Cells.Item(7, 19).Value = 0.003077940003767202
$ws.Cells.Item(7, 20).Value = 0.003077940003767201

$ws.Cells.Item(8, 9).Value = 0.04786922362394307
$ws.Cells.Item(8, 10).Value = 0.04786922362394307
$ws.Cells.Item(8, 14).Value = 49.84314599999999
$ws.Cells.Item(8, 15).Value = 0.07558164617031606
$ws.Cells.Item(8, 16).Value = 0.07558164617031604
$ws.Cells.Item(8, 17).Value = 309.9523060071386
$ws.Cells.Item(8, 18).Value = 2789.570754064247
$ws.Cells.Item(8, 19).Value = 0.0036180347223926
$ws.Cells.Item(8, 20).Value = 0.003618034722392599

$ws.Cells.Item(9, 9).Value = 0.04786922362394307
$ws.Cells.Item(9, 10).Value = 0.04786922362394307
$ws.Cells.Item(9, 13).Value = 87.857732
$ws.Cells.Item(9, 14).Value = 263.573196
$ws.Cells.Item(9, 15).Value = 0.3996797481453391
$ws.Cells.Item(9, 16).Value = 0.399679748145339
$ws.Cells.Item(9, 17).Value = 1639.044210850406
$ws.Cells.Item(9, 18).Value = 14751.39789765365
$ws.Cells.Item(9, 19).Value = 0.01913235924193048
$ws.Cells.Item(9, 20).Value = 0.01913235924193048

$ws.Cells.Item(10, 9).Value = 0.04786922362394307
$ws.Cells.Item(10, 10).Value = 0.04786922362394307
$ws.Cells.Item(10, 13).Value = 5.425038666666667
$ws.Cells.Item(10, 14).Value = 16.275116
$ws.Cells.Item(10, 15).Value = 0.02467942250059516
$ws.Cells.Item(10, 16).Value = 0.02467942250059516
$ws.Cells.Item(10, 17).Value = 101.2076913189564
$ws.Cells.Item(10, 18).Value = 910.869221870608
$ws.Cells.Item(10, 19).Value = 0.001181384794590762
$ws.Cells.Item(10, 20).Value = 0.001181384794590762

$ws.Cells.Item(11, 9).Value = 0.04786922362394307
$ws.Cells.Item(11, 10).Value = 0.04786922362394307
$ws.Cells.Item(11, 13).Value = 95.78895966666666
$ws.Cells.Item(11, 14).Value = 287.366879
$ws.Cells.Item(11, 15).Value = 0.4357602501584878
$ws.Cells.Item(11, 16).Value = 0.4357602501584877
$ws.Cells.Item(11, 17).Value = 1787.006518732272
$ws.Cells.Item(11, 18).Value = 16083.05866859045
$ws.Cells.Item(11, 19).Value = 0.02085950486126203
$ws.Cells.Item(11, 20).Value = 0.02085950486126202

$ws.Cells.Item(12, 7).Value = 171.0598806666667
$ws.Cells.Item(12, 8).Value = 513.1796420000001
$ws.Cells.Item(12, 9).Value = 0.4389285884413335
$ws.Cells.Item(12, 10).Value = 0.4389285884413335
$ws.Cells.Item(12, 13).Value = 14.13421233333333
$ws.Cells.Item(12, 14).Value = 42.402637
$ws.Cells.Item(12, 15).Value = 0.06429893302526193
$ws.Cells.Item(12, 16).Value = 0.06429893302526192
$ws.Cells.Item(12, 17).Value = 2417.796675057328
$ws.Cells.Item(12, 18).Value = 21760.17007551595
$ws.Cells.Item(12, 19).Value = 0.02822263991106206
$ws.Cells.Item(12, 20).Value = 0.02822263991106206

$ws.Cells.Item(13, 7).Value = 171.0598806666667
$ws.Cells.Item(13, 8).Value = 513.1796420000001
$ws.Cells.Item(13, 9).Value = 0.4389285884413335
$ws.Cells.Item(13, 10).Value = 0.4389285884413335
$ws.Cells.Item(13, 14).Value = 49.84314599999999
$ws.Cells.Item(13, 15).Value = 0.07558164617031606
$ws.Cells.Item(13, 16).Value = 0.07558164617031604
$ws.Cells.Item(13, 17).Value = 2842.054202270414
$ws.Cells.Item(13, 18).Value = 25578.48782043373
$ws.Cells.Item(13, 19).Value = 0.03317494526560914
$ws.Cells.Item(13, 20).Value = 0.03317494526560914

$ws.Cells.Item(14, 7).Value = 171.0598806666667
$ws.Cells.Item(14, 8).Value = 513.1796420000001
$ws.Cells.Item(14, 9).Value = 0.4389285884413335
$ws.Cells.Item(14, 10).Value = 0.4389285884413335
$ws.Cells.Item(14, 13).Value = 87.857732
$ws.Cells.Item(14, 14).Value = 263.573196
$ws.Cells.Item(14, 15).Value = 0.3996797481453391
$ws.Cells.Item(14, 16).Value = 0.399679748145339
$ws.Cells.Item(14, 17).Value = 15028.93315156398
$ws.Cells.Item(14, 18).Value = 135260.3983640758
$ws.Cells.Item(14, 19).Value = 0.1754308676820214
$ws.Cells.Item(14, 20).Value = 0.1754308676820213

$ws.Cells.Item(15, 7).Value = 171.0598806666667
$ws.Cells.Item(15, 8).Value = 513.1796420000001
$ws.Cells.Item(15, 9).Value = 0.4389285884413335
$ws.Cells.Item(15, 10).Value = 0.4389285884413335
$ws.Cells.Item(15, 13).Value = 5.425038666666667
$ws.Cells.Item(15, 14).Value = 16.275116
$ws.Cells.Item(15, 15).Value = 0.02467942250059516
$ws.Cells.Item(15, 16).Value = 0.02467942250059516
$ws.Cells.Item(15, 17).Value = 928.0064669320525
$ws.Cells.Item(15, 18).Value = 8352.058202388473
$ws.Cells.Item(15, 19).Value = 0.01083250408173352
$ws.Cells.Item(15, 20).Value = 0.01083250408173352

$ws.Cells.Item(16, 7).Value = 171.0598806666667
$ws.Cells.Item(16, 8).Value = 513.1796420000001
$ws.Cells.Item(16, 9).Value = 0.4389285884413335
$ws.Cells.Item(16, 10).Value = 0.4389285884413335
$ws.Cells.Item(16, 13).Value = 95.78895966666666
$ws.Cells.Item(16, 14).Value = 287.366879
$ws.Cells.Item(16, 15).Value = 0.4357602501584878
$ws.Cells.Item(16, 16).Value = 0.4357602501584877
$ws.Cells.Item(16, 17).Value = 16385.64800976415
$ws.Cells.Item(16, 18).Value = 147470.8320878773
$ws.Cells.Item(16, 19).Value = 0.1912676315009074
$ws.Cells.Item(16, 20).Value = 0.1912676315009074

$ws.Cells.Item(17, 7).Value = 12.628047
$ws.Cells.Item(17, 8).Value = 37.884141
$ws.Cells.Item(17, 9).Value = 0.0324027517316099
$ws.Cells.Item(17, 10).Value = 0.0324027517316099
$ws.Cells.Item(17, 13).Value = 14.13421233333333
$ws.Cells.Item(17, 14).Value = 42.402637
$ws.Cells.Item(17, 15).Value = 0.06429893302526193
$ws.Cells.Item(17, 16).Value = 0.06429893302526192
$ws.Cells.Item(17, 17).Value = 178.487497653313
$ws.Cells.Item(17, 18).Value = 1606.387478879817
$ws.Cells.Item(17, 19).Value = 0.002083462363424975
$ws.Cells.Item(17, 20).Value = 0.002083462363424974

$ws.Cells.Item(18, 7).Value = 12.628047
$ws.Cells.Item(18, 8).Value = 37.884141
$ws.Cells.Item(18, 9).Value = 0.0324027517316099
$ws.Cells.Item(18, 10).Value = 0.0324027517316099
$ws.Cells.Item(18, 14).Value = 49.84314599999999
$ws.Cells.Item(18, 15).Value = 0.07558164617031606
$ws.Cells.Item(18, 16).Value = 0.07558164617031604
$ws.Cells.Item(18, 17).Value = 209.807196771954
$ws.Cells.Item(18, 18).Value = 1888.264770947586
$ws.Cells.Item(18, 19).Value = 0.002449053316323135
$ws.Cells.Item(18, 20).Value = 0.002449053316323135

$ws.Cells.Item(19, 7).Value = 12.628047
$ws.Cells.Item(19, 8).Value = 37.884141
$ws.Cells.Item(19, 9).Value = 0.0324027517316099
$ws.Cells.Item(19, 10).Value = 0.0324027517316099
$ws.Cells.Item(19, 13).Value = 87.857732
$ws.Cells.Item(19, 14).Value = 263.573196
$ws.Cells.Item(19, 15).Value = 0.3996797481453391
$ws.Cells.Item(19, 16).Value = 0.399679748145339
$ws.Cells.Item(19, 17).Value = 1109.471569009404
$ws.Cells.Item(19, 18).Value = 9985.244121084635
$ws.Cells.Item(19, 19).Value = 0.01295072365130579
$ws.Cells.Item(19, 20).Value = 0.01295072365130579

$ws.Cells.Item(20, 7).Value = 12.628047
$ws.Cells.Item(20, 8).Value = 37.884141
$ws.Cells.Item(20, 9).Value = 0.0324027517316099
$ws.Cells.Item(20, 10).Value = 0.0324027517316099
$ws.Cells.Item(20, 13).Value = 5.425038666666667
$ws.Cells.Item(20, 14).Value = 16.275116
$ws.Cells.Item(20, 15).Value = 0.02467942250059516
$ws.Cells.Item(20, 16).Value = 0.02467942250059516
$ws.Cells.Item(20, 17).Value = 68.507643259484
$ws.Cells.Item(20, 18).Value = 616.568789335356
$ws.Cells.Item(20, 19).Value = 0.0007996812001662922
$ws.Cells.Item(20, 20).Value = 0.0007996812001662921

$ws.Cells.Item(21, 7).Value = 12.628047
$ws.Cells.Item(21, 8).Value = 37.884141
$ws.Cells.Item(21, 9).Value = 0.0324027517316099
$ws.Cells.Item(21, 10).Value = 0.0324027517316099
$ws.Cells.Item(21, 13).Value = 95.78895966666666
$ws.Cells.Item(21, 14).Value = 287.366879
$ws.Cells.Item(21, 15).Value = 0.4357602501584878
$ws.Cells.Item(21, 16).Value = 0.4357602501584877
$ws.Cells.Item(21, 17).Value = 1209.627484751771
$ws.Cells.Item(21, 18).Value = 10886.64736276594
$ws.Cells.Item(21, 19).Value = 0.0141198312003897
$ws.Cells.Item(21, 20).Value = 0.0141198312003897

$ws.Cells.Item(22, 7).Value = 139.5154473333333
$ws.Cells.Item(22, 8).Value = 418.546342
$ws.Cells.Item(22, 9).Value = 0.3579876130225438
$ws.Cells.Item(22, 10).Value = 0.3579876130225438
$ws.Cells.Item(22, 13).Value = 14.13421233333333
$ws.Cells.Item(22, 14).Value = 42.402637
$ws.Cells.Item(22, 15).Value = 0.06429893302526193
$ws.Cells.Item(22, 16).Value = 0.06429893302526192
$ws.Cells.Item(22, 17).Value = 1971.940956389317
$ws.Cells.Item(22, 18).Value = 17747.46860750385
$ws.Cells.Item(22, 19).Value = 0.02301822155360993
$ws.Cells.Item(22, 20).Value = 0.02301822155360993

$ws.Cells.Item(23, 7).Value = 139.5154473333333
$ws.Cells.Item(23, 8).Value = 418.546342
$ws.Cells.Item(23, 9).Value = 0.3579876130225438
$ws.Cells.Item(23, 10).Value = 0.3579876130225438
$ws.Cells.Item(23, 14).Value = 49.84314599999999
$ws.Cells.Item(23, 15).Value = 0.07558164617031606
$ws.Cells.Item(23, 16).Value = 0.07558164617031604
$ws.Cells.Item(23, 17).Value = 2317.962936896881
$ws.Cells.Item(23, 18).Value = 20861.66643207193
$ws.Cells.Item(23, 19).Value = 0.02705729310082594
$ws.Cells.Item(23, 20).Value = 0.02705729310082593

$ws.Cells.Item(24, 7).Value = 139.5154473333333
$ws.Cells.Item(24, 8).Value = 418.546342
$ws.Cells.Item(24, 9).Value = 0.3579876130225438
$ws.Cells.Item(24, 10).Value = 0.3579876130225438
$ws.Cells.Item(24, 13).Value = 87.857732
$ws.Cells.Item(24, 14).Value = 263.573196
$ws.Cells.Item(24, 15).Value = 0.3996797481453391
$ws.Cells.Item(24, 16).Value = 0.399679748145339
$ws.Cells.Item(24, 17).Value = 12257.51078167211
$ws.Cells.Item(24, 18).Value = 110317.597035049
$ws.Cells.Item(24, 19).Value = 0.1430803990120015
$ws.Cells.Item(24, 20).Value = 0.1430803990120014

$ws.Cells.Item(25, 7).Value = 139.5154473333333
$ws.Cells.Item(25, 8).Value = 418.546342
$ws.Cells.Item(25, 9).Value = 0.3579876130225438
$ws.Cells.Item(25, 10).Value = 0.3579876130225438
$ws.Cells.Item(25, 13).Value = 5.425038666666667
$ws.Cells.Item(25, 14).Value = 16.275116
$ws.Cells.Item(25, 15).Value = 0.02467942250059516
$ws.Cells.Item(25, 16).Value = 0.02467942250059516
$ws.Cells.Item(25, 17).Value = 756.8766963806302
$ws.Cells.Item(25, 18).Value = 6811.890267425672
$ws.Cells.Item(25, 19).Value = 0.008834927551762921
$ws.Cells.Item(25, 20).Value = 0.008834927551762921

$ws.Cells.Item(26, 7).Value = 139.5154473333333
$ws.Cells.Item(26, 8).Value = 418.546342
$ws.Cells.Item(26, 9).Value = 0.3579876130225438
$ws.Cells.Item(26, 10).Value = 0.3579876130225438
$ws.Cells.Item(26, 13).Value = 95.78895966666666
$ws.Cells.Item(26, 14).Value = 287.366879
$ws.Cells.Item(26, 15).Value = 0.4357602501584878
$ws.Cells.Item(26, 16).Value = 0.4357602501584877
$ws.Cells.Item(26, 17).Value = 13364.03955748962
$ws.Cells.Item(26, 18).Value = 120276.3560174066
$ws.Cells.Item(26, 19).Value = 0.1559967718043436
$ws.Cells.Item(26, 20).Value = 0.1559967718043436
